$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared-string table order for the "Requisitos" rows (B26:C29) changed:
# before: LOM3231, LOM3206, LOM3215, LOM3234
# after:  LOM3206, LOM3215, LOM3231, LOM3234
# So the visible text of rows 26-28 must be updated accordingly (row 29 unchanged).

$ws.Range("B26").Value = "LOM3206 -  Eletrônica  (Requisito)`n"
$ws.Range("C26").Value = "LOM3206 -  Eletrônica  (Requisito)`n"

$ws.Range("B27").Value = "LOM3215 -  Física do Estado Sólido  (Requisito)`n"
$ws.Range("C27").Value = "LOM3215 -  Física do Estado Sólido  (Requisito)`n"

$ws.Range("B28").Value = "LOM3231 -  Métodos Experimentais da Física IV  (Indicação de Conjunto)`n"
$ws.Range("C28").Value = "LOM3231 -  Métodos Experimentais da Física IV  (Indicação de Conjunto)`n"
